# Commit: "fix a couple of styles and add meta tags for yandex search engine"
#
# This script:
#  1. Fixes a few product-name strings on both sheets (comma -> period in
#     decimal sizes, and adds a "*" marker to some "нерж\нерж"/"нерж\оц"
#     fitting names on the TERMO sheet).
#  2. Switches the active sheet from sheet 2 (TERMO) to sheet 1, changes the
#     zoom level on sheet 2, and updates the selected cell on each sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1. Text fixes on sheet 1 ("1.0_aisi_321")
# ---------------------------------------------------------------------
$ws1.Range("A3").Value = "Труба 0.5м"
$ws1.Range("A4").Value = "Труба 0.3м"

# ---------------------------------------------------------------------
# 2. Text fixes on sheet 2 ("ТЕРМО_1.0_aisi_321")
# ---------------------------------------------------------------------
$ws2.Range("A4").Value = "Труба 0.5м нерж\нерж"
$ws2.Range("A5").Value = "Труба 0.5м нерж\оц"
$ws2.Range("A6").Value = "Труба 0.25м нерж\нерж"
$ws2.Range("A7").Value = "Труба 0.25м нерж\оц"

$ws2.Range("A8").Value = "Тройник 87* нерж\нерж"
$ws2.Range("A9").Value = "Тройник 87* нерж\оц"
$ws2.Range("A10").Value = "Тройник 45* нерж\нерж"
$ws2.Range("A11").Value = "Тройник 45* нерж\оц"
$ws2.Range("A12").Value = "Колено 90* нерж\нерж"
$ws2.Range("A13").Value = "Колено 90* нерж\оц"
$ws2.Range("A14").Value = "Колено 45* нерж\нерж"
$ws2.Range("A15").Value = "Колено 45* нерж\оц"

# ---------------------------------------------------------------------
# 3. View/selection changes
# ---------------------------------------------------------------------
# Sheet 2 loses the "active" flag, zoom changes 167 -> 125, selection moves
# to A11.
$ws2.Activate()
$excel.ActiveWindow.Zoom = 125
$ws2.Range("A11").Select()

# Sheet 1 becomes the active tab, keeps its zoom (135), selection moves to
# A6.
$ws1.Activate()
$ws1.Range("A6").Select()
